$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02066213248981796
$ws.Range("D2").Value = 0.2489056027130232
$ws.Range("E2").Value = 0.1807104754608204
$ws.Range("F2").Value = 0.9993407371113818
$ws.Range("G2").Value = 0.4687505247309076
$ws.Range("H2").Value = 0.6144515826206884
$ws.Range("I2").Value = 0.4704568079904092
$ws.Range("J2").Value = 0.1769199659630942
$ws.Range("K2").Value = 1.058472516520965
$ws.Range("N2").Value = 1.059454565646277
$ws.Range("O2").Value = 2.124451967642727

$ws.Range("C3").Value = 0.018137701825367
$ws.Range("D3").Value = 0.2429044160737277
$ws.Range("E3").Value = 0.1764126774282246
$ws.Range("F3").Value = 0.9978244949901338
$ws.Range("G3").Value = 0.4683098009098998
$ws.Range("H3").Value = 0.6181737008242152
$ws.Range("I3").Value = 0.4697686460050079
$ws.Range("J3").Value = 0.1727414055249383
$ws.Range("K3").Value = 0.9341490024015116
$ws.Range("N3").Value = 1.052374946025523
$ws.Range("O3").Value = 2.131032441661816

$ws.Range("C4").Value = 0.01658027940194984
$ws.Range("D4").Value = 0.2393160958367559
$ws.Range("E4").Value = 0.1738621112502301
$ws.Range("F4").Value = 0.9975097583419128
$ws.Range("G4").Value = 0.4684100999027407
$ws.Range("H4").Value = 0.620765225742872
$ws.Range("I4").Value = 0.4696709163160691
$ws.Range("J4").Value = 0.1702772835565511
$ws.Range("K4").Value = 0.8576297687980343
$ws.Range("N4").Value = 1.048394256316783
$ws.Range("O4").Value = 2.136490235872174

$ws.Range("C5").Value = 0.01594378141360409
$ws.Range("D5").Value = 0.2378782072684373
$ws.Range("E5").Value = 0.1728449969052086
$ws.Range("F5").Value = 0.9975364408591645
$ws.Range("G5").Value = 0.468544088841945
$ws.Range("H5").Value = 0.6218982858957816
$ws.Range("I5").Value = 0.4697127112641617
$ws.Range("J5").Value = 0.1692986832509149
$ws.Range("K5").Value = 0.82640366459799
$ws.Range("N5").Value = 1.046864633459833
$ws.Range("O5").Value = 2.139070440033748

$ws.Range("C6").Value = 0.01583798135989412
$ws.Range("D6").Value = 0.2376409234159098
$ws.Range("E6").Value = 0.1726774521417482
$ws.Range("F6").Value = 0.9975502289387563
$ws.Range("G6").Value = 0.4685719571715623
$ws.Range("H6").Value = 0.6220910802591391
$ws.Range("I6").Value = 0.4697245793792391
$ws.Range("J6").Value = 0.169137731088
$ws.Range("K6").Value = 0.8212160132355848
$ws.Range("N6").Value = 1.046616243330106
$ws.Range("O6").Value = 2.139520380354952

$ws.Range("C7").Value = 0.01657170276379105
$ws.Range("D7").Value = 0.2392966050743439
$ws.Range("E7").Value = 0.1738483038740988
$ws.Range("F7").Value = 0.9975094908703355
$ws.Range("G7").Value = 0.4684115301010436
$ws.Range("H7").Value = 0.6207801948298481
$ws.Range("I7").Value = 0.4696711495710737
$ws.Range("J7").Value = 0.1702639823371399
$ws.Range("K7").Value = 0.8572088171878534
$ws.Range("N7").Value = 1.048373251993638
$ws.Range("O7").Value = 2.136523591933951

$ws.Range("C8").Value = 0.01979326322836528
$ws.Range("D8").Value = 0.2468164560051775
$ws.Range("E8").Value = 0.1792103016074478
$ws.Range("F8").Value = 0.9986900210684553
$ws.Range("G8").Value = 0.4685214895400236
$ws.Range("H8").Value = 0.6156714397388399
$ws.Range("I8").Value = 0.4701520823823628
$ws.Range("J8").Value = 0.1754581350140541
$ws.Range("K8").Value = 1.015645375896469
$ws.Range("N8").Value = 1.056937766502884
$ws.Range("O8").Value = 2.126426559692931

$ws.Range("C9").Value = 0.02605097711658289
$ws.Range("D9").Value = 0.2623231413283946
$ws.Range("E9").Value = 0.1904239312949585
$ws.Range("F9").Value = 1.005897141346566
$ws.Range("G9").Value = 0.4716876839307957
$ws.Range("H9").Value = 0.6080818616907493
$ws.Range("I9").Value = 0.4736756179084409
$ws.Range("J9").Value = 0.1864493893032346
$ws.Range("K9").Value = 1.324787691395954
$ws.Range("N9").Value = 1.076620825952816
$ws.Range("O9").Value = 2.117888131863737

$ws.Range("C10").Value = 0.03061123769862206
$ws.Range("D10").Value = 0.2741742373153784
$ws.Range("E10").Value = 0.199087164292429
$ws.Range("F10").Value = 1.014180410066686
$ws.Range("G10").Value = 0.4758242884286261
$ws.Range("H10").Value = 0.6039863838863369
$ws.Range("I10").Value = 0.4778431175225819
$ws.Range("J10").Value = 0.1950167628396287
$ws.Range("K10").Value = 1.550868701455272
$ws.Range("N10").Value = 1.09282170123717
$ws.Range("O10").Value = 2.118505595069053

$ws.Range("C11").Value = 0.03267756551424839
$ws.Range("D11").Value = 0.2796641212814563
$ws.Range("E11").Value = 0.2031202512779302
$ws.Range("F11").Value = 1.018598939261992
$ws.Range("G11").Value = 0.4781018108287469
$ws.Range("H11").Value = 0.6024447785457596
$ws.Range("I11").Value = 0.4800831139910997
$ws.Range("J11").Value = 0.199021420112004
$ws.Range("K11").Value = 1.65347183489223
$ws.Range("N11").Value = 1.100565526040256
$ws.Range("O11").Value = 2.120288336146018

$ws.Range("C12").Value = 0.03345883362112545
$ws.Range("D12").Value = 0.2817570836780874
$ws.Range("E12").Value = 0.2046606816958274
$ws.Range("F12").Value = 1.020365717638924
$ws.Range("G12").Value = 0.4790213417126381
$ws.Range("H12").Value = 0.6019072351617183
$ws.Range("I12").Value = 0.4809809145104538
$ws.Range("J12").Value = 0.2005533095229453
$ws.Range("K12").Value = 1.692287987303473
$ws.Range("N12").Value = 1.1035513002643
$ws.Range("O12").Value = 2.121179769498895

$ws.Range("C13").Value = 0.03329062768131053
$ws.Range("D13").Value = 0.2813057035631061
$ws.Range("E13").Value = 0.2043283368912157
$ws.Range("F13").Value = 1.019981047706594
$ws.Range("G13").Value = 0.4788207627863841
$ws.Range("H13").Value = 0.6020209486122354
$ws.Range("I13").Value = 0.4807853517005327
$ws.Range("J13").Value = 0.200222704561682
$ws.Range("K13").Value = 1.683929940159715
$ws.Range("N13").Value = 1.10290589396142
$ws.Range("O13").Value = 2.120978155174157

$ws.Range("C14").Value = 0.03274186525457878
$ws.Range("D14").Value = 0.2798360295222864
$ws.Range("E14").Value = 0.203246719463877
$ws.Range("F14").Value = 1.018742417494153
$ws.Range("G14").Value = 0.4781763161840189
$ws.Range("H14").Value = 0.6023996279981247
$ws.Range("I14").Value = 0.4801559828784008
$ws.Range("J14").Value = 0.1991471407963132
$ws.Range("K14").Value = 1.656666028034863
$ws.Range("N14").Value = 1.100810100533167
$ws.Range("O14").Value = 2.120357336785162

$ws.Range("C15").Value = 0.03240557448447134
$ws.Range("D15").Value = 0.2789376389532237
$ws.Range("E15").Value = 0.2025859131814158
$ws.Range("F15").Value = 1.017995908598394
$ws.Range("G15").Value = 0.4777890129698932
$ws.Range("H15").Value = 0.6026376006046377
$ws.Range("I15").Value = 0.4797769328524453
$ws.Range("J15").Value = 0.1984903335303727
$ws.Range("K15").Value = 1.63996114366978
$ws.Range("N15").Value = 1.099533302702
$ws.Range("O15").Value = 2.120005253808841

$ws.Range("C16").Value = 0.03047603140210242
$ws.Range("D16").Value = 0.273817435911738
$ws.Range("E16").Value = 0.1988254413877257
$ws.Range("F16").Value = 1.013904741081845
$ws.Range("G16").Value = 0.4756834254914253
$ws.Range("H16").Value = 0.6040936002277704
$ws.Range("I16").Value = 0.4777036602584346
$ws.Range("J16").Value = 0.1947572069026506
$ws.Range("K16").Value = 1.544158248926351
$ws.Range("N16").Value = 1.092323113159893
$ws.Range("O16").Value = 2.11841934272303

$ws.Range("C17").Value = 0.02929020688073081
$ws.Range("D17").Value = 0.2707015532143515
$ws.Range("E17").Value = 0.1965420690149742
$ws.Range("F17").Value = 1.011561566210645
$ws.Range("G17").Value = 0.4744932028472988
$ws.Range("H17").Value = 0.6050691431544664
$ws.Range("I17").Value = 0.4765199774362259
$ws.Range("J17").Value = 0.1924945272612746
$ws.Range("K17").Value = 1.485322383023458
$ws.Range("N17").Value = 1.087995375945752
$ws.Range("O17").Value = 2.117831372188334

$ws.Range("C18").Value = 0.02860738608826807
$ws.Range("D18").Value = 0.268918684077903
$ws.Range("E18").Value = 0.1952374098630187
$ws.Range("F18").Value = 1.010275052743822
$ws.Range("G18").Value = 0.4738458588571604
$ws.Range("H18").Value = 0.6056605055411239
$ws.Range("I18").Value = 0.4758715461264842
$ws.Range("J18").Value = 0.1912031960600871
$ws.Range("K18").Value = 1.451458901668786
$ws.Range("N18").Value = 1.085541400016055
$ws.Range("O18").Value = 2.117634525978787

$ws.Range("C19").Value = 0.02837606414200877
$ws.Range("D19").Value = 0.268316637464352
$ws.Range("E19").Value = 0.194797166714558
$ws.Range("F19").Value = 1.009849974331175
$ws.Range("G19").Value = 0.4736330701268656
$ws.Range("H19").Value = 0.6058659270037907
$ws.Range("I19").Value = 0.4756575595959589
$ws.Range("J19").Value = 0.1907677091606104
$ws.Range("K19").Value = 1.439989498115608
$ws.Range("N19").Value = 1.084716588605289
$ws.Range("O19").Value = 2.117592140599754

$ws.Range("C20").Value = 0.02941651942700219
$ws.Range("D20").Value = 0.2710322823070754
$ws.Range("E20").Value = 0.1967842403887587
$ws.Range("F20").Value = 1.011804664923915
$ws.Range("G20").Value = 0.4746160485614155
$ws.Range("H20").Value = 0.604962163533969
$ws.Range("I20").Value = 0.4766426295601889
$ws.Range("J20").Value = 0.1927343479627694
$ws.Range("K20").Value = 1.491587924110377
$ws.Range("N20").Value = 1.088452428466056
$ws.Range("O20").Value = 2.117879331821513

$ws.Range("C21").Value = 0.03290308305474809
$ws.Range("D21").Value = 0.2802673278720818
$ws.Range("E21").Value = 0.2035640592971291
$ws.Range("F21").Value = 1.019103693189763
$ws.Range("G21").Value = 0.47836405519854
$ws.Range("H21").Value = 0.602287145996911
$ws.Range("I21").Value = 0.4803394980953826
$ws.Range("J21").Value = 0.1994626418085375
$ws.Range("K21").Value = 1.664675136303401
$ws.Range("N21").Value = 1.101424241116703
$ws.Range("O21").Value = 2.120533811534642

$ws.Range("C22").Value = 0.0351747129188027
$ws.Range("D22").Value = 0.2863848662451005
$ws.Range("E22").Value = 0.2080719033137086
$ws.Range("F22").Value = 1.024419490500691
$ws.Range("G22").Value = 0.4811463942624528
$ws.Range("H22").Value = 0.600808328386222
$ws.Range("I22").Value = 0.4830445350365551
$ws.Range("J22").Value = 0.2039497980152163
$ws.Range("K22").Value = 1.777578281817568
$ws.Range("N22").Value = 1.110212879136654
$ws.Range("O22").Value = 2.123529884669921

$ws.Range("C23").Value = 0.0339629562704431
$ws.Range("D23").Value = 0.2831123728498142
$ws.Range("E23").Value = 0.2056589716331274
$ws.Range("F23").Value = 1.021532424822283
$ws.Range("G23").Value = 0.4796308996258745
$ws.Range("H23").Value = 0.601572943753041
$ws.Range("I23").Value = 0.4815743469588369
$ws.Range("J23").Value = 0.2015467055982469
$ws.Range("K23").Value = 1.717340673604497
$ws.Range("N23").Value = 1.105493915601485
$ws.Range("O23").Value = 2.121815296833404

$ws.Range("C24").Value = 0.02935941687071875
$ws.Range("D24").Value = 0.2708827332192101
$ws.Range("E24").Value = 0.1966747295322548
$ws.Range("F24").Value = 1.011694571206377
$ws.Range("G24").Value = 0.4745603949902204
$ws.Range("H24").Value = 0.6050104339839208
$ws.Range("I24").Value = 0.4765870786017175
$ws.Range("J24").Value = 0.1926258953899378
$ws.Range("K24").Value = 1.48875539102238
$ws.Range("N24").Value = 1.088245688776595
$ws.Range("O24").Value = 2.117857209513744

$ws.Range("C25").Value = 0.0243645848930214
$ws.Range("D25").Value = 0.2580473112721364
$ws.Range("E25").Value = 0.1873156840441155
$ws.Range("F25").Value = 1.003423082954079
$ws.Range("G25").Value = 0.4705140565178709
$ws.Range("H25").Value = 0.6098750375550424
$ws.Range("I25").Value = 0.4724455672794221
$ws.Range("J25").Value = 0.1833896184022592
$ws.Range("K25").Value = 1.241332960551915
$ws.Range("N25").Value = 1.070989038282107
$ws.Range("O25").Value = 2.118989739447272
